$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.379.65"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "2.095.17"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'251.94"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").Value = "'0.668"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'54.44"
$ws.Range("E8").Value = "  +20.57%  "
$ws.Range("D9").Value = "'62.22"
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("D10").Value = "'0.381"
$ws.Range("E10").Value = "  +4.79%  "
$ws.Range("E11").Value = "  +4.29%  "
$ws.Range("E12").Value = "  +7.59%  "
$ws.Range("D13").Value = "'15.08"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "2.395.74"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "2.092.55"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "'5.23"
$ws.Range("E17").Value = "  +6.69%  "
$ws.Range("D18").Value = "37.351.00"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").Value = "'73.18"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").Value = "'14.47"
$ws.Range("E20").Value = "  +14.65%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +4.58%  "
$ws.Range("D22").Value = "'241.22"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").Value = "'5.25"
$ws.Range("E23").Value = "  +6.67%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "'171.76"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("D28").Value = "'20.94"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").Value = "'2.04"
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +26.51%  "
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("D34").Value = "'0.0623"
$ws.Range("E34").Value = "  +7.10%  "
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +6.18%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.25"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("B40").Value = "FTXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D40").Value = "'5.01"
$ws.Range("E40").Value = "  +144.99%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.35"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("E42").Value = "  +12.89%  "
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("D44").Value = "'1.18"
$ws.Range("E44").Value = "  +4.87%  "
$ws.Range("D45").Value = "'0.0973"
$ws.Range("E45").Value = "  +19.37%  "
$ws.Range("D46").Value = "'99.26"
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("D47").Value = "'2.81"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").Value = "1.331.76"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = "  +4.95%  "
$ws.Range("E50").Value = "  +7.41%  "
$ws.Range("D51").Value = "'6.94"
$ws.Range("E51").Value = "  +13.84%  "
